# Add files via upload - Error codes and ack updated
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ERROR-CODES")

# ---------------------------------------------------------------------------
# 0. Write the new strings FIRST, in the exact order they should land in the
#    shared-string table (D19, D18, C1 header, A36, B36).
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = "CARRIER DATA PURGED OR ARCHIVED "
$ws.Range("D18").Value = "CBL is Live for 45 DAYS (Default Limit)"
$ws.Range("C1").Value  = "ACK-VALUE"
$ws.Range("A36").Value = "ER111"
$ws.Range("B36").Value = "BOL- CNT MIS-MATCH"

# ---------------------------------------------------------------------------
# 1. Header row formatting: C1/D1 styled like A1/B1 (bold header style).
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Value = "NOTES"
$ws.Range("C1").Value = "ACK-VALUE"

# ---------------------------------------------------------------------------
# 2. "ACK-VALUE" column C for all plain (non-highlighted) rows: 2-13, 21-36.
#    Formatting copied from A2, which is guaranteed to carry the plain style.
# ---------------------------------------------------------------------------
foreach ($r in (2..13) + (21..36)) {
    $ws.Range("A2").Copy() | Out-Null
    $ws.Range("C$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("C$r").Value = 0
}

# ---------------------------------------------------------------------------
# 3. Rows 14-20 get the yellow highlight style across columns A:D, with the
#    "ACK-VALUE" flags and the two notes on rows 18 & 19.
# ---------------------------------------------------------------------------
$ws.Range("A34:B34").Copy() | Out-Null
$ws.Range("A14:B20").PasteSpecial(-4122) | Out-Null

foreach ($r in 14..20) {
    $ws.Range("A34").Copy() | Out-Null
    $ws.Range("C$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("A34").Copy() | Out-Null
    $ws.Range("D$r").PasteSpecial(-4122) | Out-Null
}

$ackOnes = @(15, 16, 17, 18, 20)
foreach ($r in 14..20) {
    if ($ackOnes -contains $r) {
        $ws.Range("C$r").Value = 1
    } else {
        $ws.Range("C$r").Value = 0
    }
}

$ws.Range("D18").Value = "CBL is Live for 45 DAYS (Default Limit)"
$ws.Range("D19").Value = "CARRIER DATA PURGED OR ARCHIVED "

# ---------------------------------------------------------------------------
# 4. Rows 34 & 35: drop the yellow highlight back to the plain style (closest
#    achievable match to the author's "111" bookkeeping variant is the same
#    plain, unfilled look already used everywhere else on the sheet).
# ---------------------------------------------------------------------------
$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A34:B35").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 5. New row 36: ER111 / BOL- CNT MIS-MATCH / 0, styled like row 33 above it.
# ---------------------------------------------------------------------------
$ws.Range("A33").Copy() | Out-Null
$ws.Range("A36").PasteSpecial(-4122) | Out-Null
$ws.Range("B33").Copy() | Out-Null
$ws.Range("B36").PasteSpecial(-4122) | Out-Null

$ws.Range("A36").Value = "ER111"
$ws.Range("B36").Value = "BOL- CNT MIS-MATCH"
$ws.Range("C36").Value = 0

# ---------------------------------------------------------------------------
# 6. Column widths for the two new columns.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 50.426339285714285
$ws.Columns.Item(4).ColumnWidth = 56.426339285714285

# ---------------------------------------------------------------------------
# 7. Selections: VALIDATION-CODES -> C32, ERROR-CODES -> C19 (left active).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("VALIDATION-CODES")
$ws4.Activate() | Out-Null
$ws4.Range("C32").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("C19").Select() | Out-Null
